$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new header cells (AD1:AF1) ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting of the other header cells (bold, border, centered/top aligned)
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- Fill in the team record values for every data row (2-50) ---
$lastRow = 50
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 74   # AD = column 30 -> Wins
    $ws.Cells.Item($r, 31).Value = 88   # AE = column 31 -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF = column 32 -> Ties
}

$excel.CutCopyMode = 0
